# Append a new paragraph after the last paragraph of the document, carrying
# the trailing "_GoBack" bookmark along with it so that it still marks the
# very end of the document content (as it did before the edit).

$d = $word.ActiveDocument

# 1. Add a new, empty paragraph right after the current last paragraph and
#    fill it with the new sentence. A short sentinel ("###") is appended so
#    that, while we are positioning the bookmark, its end point does not sit
#    exactly on the document's final character -- doing so confuses bookmark
#    placement. We strip the sentinel off again afterwards.
$lastParagraph = $d.Paragraphs.Last
$lastParagraph.Range.InsertParagraphAfter()

$newParagraph = $d.Paragraphs.Last
$newParagraph.Range.Text = "When the interaction is strong, those species started out at disadvantage quickly dies out, the remaining combination is limited.###"

# 2. Figure out where the real text ends (i.e. just before the "###"
#    sentinel and the trailing paragraph mark) and drop a collapsed
#    "_GoBack" bookmark there -- this is where Word always keeps the
#    last-edit bookmark.
$target = $d.Paragraphs.Last
$bookmarkPos = $target.Range.End - 4
$bookmarkRange = $d.Range($bookmarkPos, $bookmarkPos)
$bookmarkRange.Bookmarks.Add("_GoBack")

# 3. Remove the sentinel text now that the bookmark is anchored in place.
$sentinelRange = $d.Range($target.Range.End - 4, $target.Range.End - 1)
$sentinelRange.Delete()
